$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F ("想去人数")
$updates = @{
    2  = 1877
    3  = 263
    4  = 251
    5  = 8385
    6  = 571
    7  = 604
    8  = 86
    10 = 9245
    11 = 2395
    12 = 8
    13 = 37
    14 = 330
    15 = 10197
    16 = 10575
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
